$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: BTC -> BTC
$ws.Cells.Item(2, 4).Value = 29363
$ws.Cells.Item(2, 5).Value = 571333545378
$ws.Cells.Item(2, 6).Value = 3303796092
$ws.Cells.Item(2, 7).Value = -0.13376

# Row 3: ETH -> ETH
$ws.Cells.Item(3, 4).Value = 1848.18
$ws.Cells.Item(3, 5).Value = 222076564062
$ws.Cells.Item(3, 6).Value = 2310450298
$ws.Cells.Item(3, 7).Value = -0.08673

# Row 4: USDT -> USDT
$ws.Cells.Item(4, 4).Value = 0.999184
$ws.Cells.Item(4, 5).Value = 83342810176
$ws.Cells.Item(4, 6).Value = 7108679688
$ws.Cells.Item(4, 7).Value = -0.00879

# Row 5: BNB -> BNB
$ws.Cells.Item(5, 4).Value = 240.58
$ws.Cells.Item(5, 5).Value = 37018923781
$ws.Cells.Item(5, 6).Value = 213453442
$ws.Cells.Item(5, 7).Value = -0.07003

# Row 6: XRP -> XRP
$ws.Cells.Item(6, 4).Value = 0.62698
$ws.Cells.Item(6, 5).Value = 33106832695
$ws.Cells.Item(6, 6).Value = 483544806
$ws.Cells.Item(6, 7).Value = -0.39903

# Row 7: USDC -> USDC
$ws.Cells.Item(7, 4).Value = 0.999383
$ws.Cells.Item(7, 5).Value = 26149288275
$ws.Cells.Item(7, 6).Value = 1179107340
$ws.Cells.Item(7, 7).Value = -0.06833

# Row 8: STETH -> STETH
$ws.Cells.Item(8, 4).Value = 1846.9
$ws.Cells.Item(8, 5).Value = 14936655272
$ws.Cells.Item(8, 6).Value = 3754359
$ws.Cells.Item(8, 7).Value = -0.17284

# Row 9: DOGE -> DOGE
$ws.Cells.Item(9, 4).Value = 0.076025
$ws.Cells.Item(9, 5).Value = 10692730010
$ws.Cells.Item(9, 6).Value = 312734566
$ws.Cells.Item(9, 7).Value = -1.02666

# Row 10: ADA -> ADA
$ws.Cells.Item(10, 4).Value = 0.290117
$ws.Cells.Item(10, 5).Value = 10171497411
$ws.Cells.Item(10, 6).Value = 107779936
$ws.Cells.Item(10, 7).Value = -1.28148

# Row 11: SOL -> SOL
$ws.Cells.Item(11, 4).Value = 24.65
$ws.Cells.Item(11, 5).Value = 10008374813
$ws.Cells.Item(11, 6).Value = 287251378
$ws.Cells.Item(11, 7).Value = 0.50687

# Row 12: TRX -> TRX
$ws.Cells.Item(12, 4).Value = 0.077405
$ws.Cells.Item(12, 5).Value = 6929030335
$ws.Cells.Item(12, 6).Value = 143140064
$ws.Cells.Item(12, 7).Value = -0.02841

# Row 13: DOT -> DOT
$ws.Cells.Item(13, 4).Value = 5.03
$ws.Cells.Item(13, 5).Value = 6354890239
$ws.Cells.Item(13, 6).Value = 66099430
$ws.Cells.Item(13, 7).Value = 0.02147

# Row 14: MATIC -> MATIC
$ws.Cells.Item(14, 4).Value = 0.678663
$ws.Cells.Item(14, 5).Value = 6325172732
$ws.Cells.Item(14, 6).Value = 93208608
$ws.Cells.Item(14, 7).Value = -0.36905

# Row 15: LTC -> SHIB
$ws.Cells.Item(15, 2).Value = 'SHIB'
$ws.Cells.Item(15, 3).Value = 'Shiba Inu'
$ws.Cells.Item(15, 4).Value = 0.00001063
$ws.Cells.Item(15, 5).Value = 6277730810
$ws.Cells.Item(15, 6).Value = 431703613
$ws.Cells.Item(15, 7).Value = -2.46102

# Row 16: SHIB -> LTC
$ws.Cells.Item(16, 2).Value = 'LTC'
$ws.Cells.Item(16, 3).Value = 'Litecoin'
$ws.Cells.Item(16, 4).Value = 82.9
$ws.Cells.Item(16, 5).Value = 6096772616
$ws.Cells.Item(16, 6).Value = 356812368
$ws.Cells.Item(16, 7).Value = -0.99549

# Row 17: WBTC -> TON
$ws.Cells.Item(17, 2).Value = 'TON'
$ws.Cells.Item(17, 3).Value = 'Toncoin'
$ws.Cells.Item(17, 4).Value = 1.43
$ws.Cells.Item(17, 5).Value = 4936033777
$ws.Cells.Item(17, 6).Value = 30373319
$ws.Cells.Item(17, 7).Value = 9.21244

# Row 18: UNI -> WBTC
$ws.Cells.Item(18, 2).Value = 'WBTC'
$ws.Cells.Item(18, 3).Value = 'Wrapped Bitcoin'
$ws.Cells.Item(18, 4).Value = 29391
$ws.Cells.Item(18, 5).Value = 4768484517
$ws.Cells.Item(18, 6).Value = 30522230
$ws.Cells.Item(18, 7).Value = -0.1164

# Row 19: BCH -> UNI
$ws.Cells.Item(19, 2).Value = 'UNI'
$ws.Cells.Item(19, 3).Value = 'Uniswap'
$ws.Cells.Item(19, 4).Value = 6.12
$ws.Cells.Item(19, 5).Value = 4617066618
$ws.Cells.Item(19, 6).Value = 53502521
$ws.Cells.Item(19, 7).Value = -0.5577

# Row 20: AVAX -> BCH
$ws.Cells.Item(20, 2).Value = 'BCH'
$ws.Cells.Item(20, 3).Value = 'Bitcoin Cash'
$ws.Cells.Item(20, 4).Value = 227.82
$ws.Cells.Item(20, 5).Value = 4438553491
$ws.Cells.Item(20, 6).Value = 98148620
$ws.Cells.Item(20, 7).Value = -0.61341

# Row 21: TON -> AVAX
$ws.Cells.Item(21, 2).Value = 'AVAX'
$ws.Cells.Item(21, 3).Value = 'Avalanche'
$ws.Cells.Item(21, 4).Value = 12.35
$ws.Cells.Item(21, 5).Value = 4245537266
$ws.Cells.Item(21, 6).Value = 82754761
$ws.Cells.Item(21, 7).Value = -1.0073

# Row 22: DAI -> LINK
$ws.Cells.Item(22, 2).Value = 'LINK'
$ws.Cells.Item(22, 3).Value = 'Chainlink'
$ws.Cells.Item(22, 4).Value = 7.5
$ws.Cells.Item(22, 5).Value = 4038541357
$ws.Cells.Item(22, 6).Value = 122835487
$ws.Cells.Item(22, 7).Value = 0.71173

# Row 23: XLM -> DAI
$ws.Cells.Item(23, 2).Value = 'DAI'
$ws.Cells.Item(23, 3).Value = 'Dai'
$ws.Cells.Item(23, 4).Value = 0.999619
$ws.Cells.Item(23, 5).Value = 3999458311
$ws.Cells.Item(23, 6).Value = 56265311
$ws.Cells.Item(23, 7).Value = -0.01704

# Row 24: LINK -> XLM
$ws.Cells.Item(24, 2).Value = 'XLM'
$ws.Cells.Item(24, 3).Value = 'Stellar'
$ws.Cells.Item(24, 4).Value = 0.138305
$ws.Cells.Item(24, 5).Value = 3786616248
$ws.Cells.Item(24, 6).Value = 58806981
$ws.Cells.Item(24, 7).Value = -0.08619

# Row 25: LEO -> LEO
$ws.Cells.Item(25, 4).Value = 4.03
$ws.Cells.Item(25, 5).Value = 3757947879
$ws.Cells.Item(25, 6).Value = 829456
$ws.Cells.Item(25, 7).Value = 1.27446

# Row 26: BUSD -> BUSD
$ws.Cells.Item(26, 4).Value = 0.999841
$ws.Cells.Item(26, 5).Value = 3371446673
$ws.Cells.Item(26, 6).Value = 953231543
$ws.Cells.Item(26, 7).Value = -0.01211

# Row 27: TUSD -> TUSD
$ws.Cells.Item(27, 4).Value = 0.999184
$ws.Cells.Item(27, 5).Value = 2970859209
$ws.Cells.Item(27, 6).Value = 61108173
$ws.Cells.Item(27, 7).Value = -0.06144

# Row 28: XMR -> XMR
$ws.Cells.Item(28, 4).Value = 159.02
$ws.Cells.Item(28, 5).Value = 2883280463
$ws.Cells.Item(28, 6).Value = 72144578
$ws.Cells.Item(28, 7).Value = 1.16487

# Row 29: OKB -> OKB
$ws.Cells.Item(29, 4).Value = 46.71
$ws.Cells.Item(29, 5).Value = 2802373037
$ws.Cells.Item(29, 6).Value = 5304234
$ws.Cells.Item(29, 7).Value = -3.00401

# Row 30: ETC -> ETC
$ws.Cells.Item(30, 4).Value = 17.66
$ws.Cells.Item(30, 5).Value = 2514932071
$ws.Cells.Item(30, 6).Value = 54094016
$ws.Cells.Item(30, 7).Value = -0.07175

# Row 31: ATOM -> ATOM
$ws.Cells.Item(31, 4).Value = 8.44
$ws.Cells.Item(31, 5).Value = 2469359258
$ws.Cells.Item(31, 6).Value = 64621810
$ws.Cells.Item(31, 7).Value = 0.67488

# Row 32: FIL -> HBAR
$ws.Cells.Item(32, 2).Value = 'HBAR'
$ws.Cells.Item(32, 3).Value = 'Hedera'
$ws.Cells.Item(32, 4).Value = 0.056172
$ws.Cells.Item(32, 5).Value = 1852562666
$ws.Cells.Item(32, 6).Value = 16324683
$ws.Cells.Item(32, 7).Value = -1.60931

# Row 33: ICP -> FIL
$ws.Cells.Item(33, 2).Value = 'FIL'
$ws.Cells.Item(33, 3).Value = 'Filecoin'
$ws.Cells.Item(33, 4).Value = 4.1
$ws.Cells.Item(33, 5).Value = 1810280167
$ws.Cells.Item(33, 6).Value = 52786199
$ws.Cells.Item(33, 7).Value = -0.41291

# Row 34: HBAR -> ICP
$ws.Cells.Item(34, 2).Value = 'ICP'
$ws.Cells.Item(34, 3).Value = 'Internet Computer'
$ws.Cells.Item(34, 4).Value = 4.06
$ws.Cells.Item(34, 5).Value = 1793068862
$ws.Cells.Item(34, 6).Value = 14078567
$ws.Cells.Item(34, 7).Value = 0.11028

# Row 35: MNT -> APT
$ws.Cells.Item(35, 2).Value = 'APT'
$ws.Cells.Item(35, 3).Value = 'Aptos'
$ws.Cells.Item(35, 4).Value = 7.21
$ws.Cells.Item(35, 5).Value = 1621145976
$ws.Cells.Item(35, 6).Value = 58254185
$ws.Cells.Item(35, 7).Value = 0.97991

# Row 36: LDO -> LDO
$ws.Cells.Item(36, 4).Value = 1.83
$ws.Cells.Item(36, 5).Value = 1608902042
$ws.Cells.Item(36, 6).Value = 64275190
$ws.Cells.Item(36, 7).Value = -1.10577

# Row 37: CRO -> MNT
$ws.Cells.Item(37, 2).Value = 'MNT'
$ws.Cells.Item(37, 3).Value = 'Mantle'
$ws.Cells.Item(37, 4).Value = 0.462679
$ws.Cells.Item(37, 5).Value = 1496606169
$ws.Cells.Item(37, 6).Value = 4522297
$ws.Cells.Item(37, 7).Value = 0.00789

# Row 38: APT -> CRO
$ws.Cells.Item(38, 2).Value = 'CRO'
$ws.Cells.Item(38, 3).Value = 'Cronos'
$ws.Cells.Item(38, 4).Value = 0.056974
$ws.Cells.Item(38, 5).Value = 1492000194
$ws.Cells.Item(38, 6).Value = 3365776
$ws.Cells.Item(38, 7).Value = -0.38157

# Row 39: QNT -> ARB
$ws.Cells.Item(39, 2).Value = 'ARB'
$ws.Cells.Item(39, 3).Value = 'Arbitrum'
$ws.Cells.Item(39, 4).Value = 1.16
$ws.Cells.Item(39, 5).Value = 1479485675
$ws.Cells.Item(39, 6).Value = 96635627
$ws.Cells.Item(39, 7).Value = 0.08483

# Row 40: ARB -> QNT
$ws.Cells.Item(40, 2).Value = 'QNT'
$ws.Cells.Item(40, 3).Value = 'Quant'
$ws.Cells.Item(40, 4).Value = 101.44
$ws.Cells.Item(40, 5).Value = 1475426292
$ws.Cells.Item(40, 6).Value = 10583241
$ws.Cells.Item(40, 7).Value = -0.05005

# Row 41: VET -> VET
$ws.Cells.Item(41, 4).Value = 0.01798175
$ws.Cells.Item(41, 5).Value = 1307156169
$ws.Cells.Item(41, 6).Value = 23665482
$ws.Cells.Item(41, 7).Value = -0.01439

# Row 42: NEAR -> NEAR
$ws.Cells.Item(42, 4).Value = 1.34
$ws.Cells.Item(42, 5).Value = 1264025411
$ws.Cells.Item(42, 6).Value = 40244806
$ws.Cells.Item(42, 7).Value = -0.57759

# Row 43: OP -> OP
$ws.Cells.Item(43, 4).Value = 1.55
$ws.Cells.Item(43, 5).Value = 1112859925
$ws.Cells.Item(43, 6).Value = 60081069
$ws.Cells.Item(43, 7).Value = -0.85655

# Row 44: MKR -> MKR
$ws.Cells.Item(44, 4).Value = 1232.86
$ws.Cells.Item(44, 5).Value = 1110384600
$ws.Cells.Item(44, 6).Value = 92594121
$ws.Cells.Item(44, 7).Value = 0.73881

# Row 45: KAS -> RETH
$ws.Cells.Item(45, 2).Value = 'RETH'
$ws.Cells.Item(45, 3).Value = 'Rocket Pool ETH'
$ws.Cells.Item(45, 4).Value = 2006.7
$ws.Cells.Item(45, 5).Value = 1005438874
$ws.Cells.Item(45, 6).Value = 5396065
$ws.Cells.Item(45, 7).Value = -0.13538

# Row 46: XDC -> GRT
$ws.Cells.Item(46, 2).Value = 'GRT'
$ws.Cells.Item(46, 3).Value = 'The Graph'
$ws.Cells.Item(46, 4).Value = 0.107169
$ws.Cells.Item(46, 5).Value = 979040189
$ws.Cells.Item(46, 6).Value = 34257496
$ws.Cells.Item(46, 7).Value = 1.27626

# Row 47: RETH -> AAVE
$ws.Cells.Item(47, 2).Value = 'AAVE'
$ws.Cells.Item(47, 3).Value = 'Aave'
$ws.Cells.Item(47, 4).Value = 65.29
$ws.Cells.Item(47, 5).Value = 948408667
$ws.Cells.Item(47, 6).Value = 56626564
$ws.Cells.Item(47, 7).Value = -1.34104

# Row 48: GRT -> KAS
$ws.Cells.Item(48, 2).Value = 'KAS'
$ws.Cells.Item(48, 3).Value = 'Kaspa'
$ws.Cells.Item(48, 4).Value = 0.04492355
$ws.Cells.Item(48, 5).Value = 911413390
$ws.Cells.Item(48, 6).Value = 17374165
$ws.Cells.Item(48, 7).Value = -2.66644

# Row 49: AAVE -> ALGO
$ws.Cells.Item(49, 2).Value = 'ALGO'
$ws.Cells.Item(49, 3).Value = 'Algorand'
$ws.Cells.Item(49, 4).Value = 0.114077
$ws.Cells.Item(49, 5).Value = 891551189
$ws.Cells.Item(49, 6).Value = 34316017
$ws.Cells.Item(49, 7).Value = 1.16763

# Row 50: ALGO -> XDC
$ws.Cells.Item(50, 2).Value = 'XDC'
$ws.Cells.Item(50, 3).Value = 'XDC Network'
$ws.Cells.Item(50, 4).Value = 0.063852
$ws.Cells.Item(50, 5).Value = 884943140
$ws.Cells.Item(50, 6).Value = 13176740
$ws.Cells.Item(50, 7).Value = -2.9785

# Row 51: SAND -> STX
$ws.Cells.Item(51, 2).Value = 'STX'
$ws.Cells.Item(51, 3).Value = 'Stacks'
$ws.Cells.Item(51, 4).Value = 0.59442
$ws.Cells.Item(51, 5).Value = 832994043
$ws.Cells.Item(51, 6).Value = 15906286
$ws.Cells.Item(51, 7).Value = 2.15945

